{"js": "// Replace the text of every arithmetic-problem cell in the worksheet's\n// table with its updated value, cell by cell, preserving each cell's\n// existing paragraph/run formatting (font, size, alignment). The table\n// is 20 rows x 5 columns of \"NN+NN=\" / \"NN-NN=\" style problems; this\n// array holds the new value for every cell in row-major (reading) order.\nconst newValues = [[\"38+24=\", \"24+48=\", \"54-46=\", \"61-32=\", \"19+16=\"], [\"16+37=\", \"76+9=\", \"38+16=\", \"29+56=\", \"34+27=\"], [\"20-15=\", \"80-8=\", \"26+6=\", \"22+19=\", \"13+49=\"], [\"97-49=\", \"16+27=\", \"39+13=\", \"86-59=\", \"92-84=\"], [\"74+9=\", \"54-19=\", \"48+33=\", \"28+26=\", \"63+19=\"], [\"17+38=\", \"56-38=\", \"85-27=\", \"93-36=\", \"31-26=\"], [\"25+58=\", \"26+69=\", \"82-78=\", \"87+4=\", \"70-14=\"], [\"36+17=\", \"81-78=\", \"42-7=\", \"5+48=\", \"19+2=\"], [\"49+7=\", \"50-37=\", \"47+16=\", \"90-23=\", \"72-65=\"], [\"9+84=\", \"46+45=\", \"90-66=\", \"76-18=\", \"9+85=\"], [\"82-66=\", \"74-36=\", \"8+46=\", \"39+6=\", \"67+24=\"], [\"70-2=\", \"42-35=\", \"9+54=\", \"93-47=\", \"56+37=\"], [\"58+23=\", \"84-59=\", \"80-11=\", \"76+9=\", \"5+77=\"], [\"70-53=\", \"54-28=\", \"23-6=\", \"82-9=\", \"9+72=\"], [\"71-45=\", \"85-7=\", \"36+19=\", \"57+4=\", \"42-29=\"], [\"60-44=\", \"46+47=\", \"53-16=\", \"42-33=\", \"41-39=\"], [\"16+6=\", \"93-25=\", \"43+29=\", \"52-49=\", \"35-18=\"], [\"49+33=\", \"87+4=\", \"32-5=\", \"21-3=\", \"62-25=\"], [\"87-59=\", \"16+78=\", \"13+8=\", \"9+19=\", \"37+14=\"], [\"62-53=\", \"56+35=\", \"40-21=\", \"86-69=\", \"84-77=\"]];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = Math.min(table.rowCount, newValues.length);\nfor (let r = 0; r < rowCount; r++) {\n  const rowValues = newValues[r];\n  for (let c = 0; c < rowValues.length; c++) {\n    const cell = table.getCell(r, c);\n    const range = cell.body.getRange();\n    range.insertText(rowValues[c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Replace the text of every arithmetic-problem cell in the worksheet's\n# table with its updated value, cell by cell, preserving each cell's\n# existing paragraph/run formatting (font, size, alignment). The table\n# is 20 rows x 5 columns of \"NN+NN=\" / \"NN-NN=\" style problems; this\n# array holds the new value for every cell in row-major (reading) order.\n$newValues = @(\n    @(\"38+24=\", \"24+48=\", \"54-46=\", \"61-32=\", \"19+16=\"),\n    @(\"16+37=\", \"76+9=\", \"38+16=\", \"29+56=\", \"34+27=\"),\n    @(\"20-15=\", \"80-8=\", \"26+6=\", \"22+19=\", \"13+49=\"),\n    @(\"97-49=\", \"16+27=\", \"39+13=\", \"86-59=\", \"92-84=\"),\n    @(\"74+9=\", \"54-19=\", \"48+33=\", \"28+26=\", \"63+19=\"),\n    @(\"17+38=\", \"56-38=\", \"85-27=\", \"93-36=\", \"31-26=\"),\n    @(\"25+58=\", \"26+69=\", \"82-78=\", \"87+4=\", \"70-14=\"),\n    @(\"36+17=\", \"81-78=\", \"42-7=\", \"5+48=\", \"19+2=\"),\n    @(\"49+7=\", \"50-37=\", \"47+16=\", \"90-23=\", \"72-65=\"),\n    @(\"9+84=\", \"46+45=\", \"90-66=\", \"76-18=\", \"9+85=\"),\n    @(\"82-66=\", \"74-36=\", \"8+46=\", \"39+6=\", \"67+24=\"),\n    @(\"70-2=\", \"42-35=\", \"9+54=\", \"93-47=\", \"56+37=\"),\n    @(\"58+23=\", \"84-59=\", \"80-11=\", \"76+9=\", \"5+77=\"),\n    @(\"70-53=\", \"54-28=\", \"23-6=\", \"82-9=\", \"9+72=\"),\n    @(\"71-45=\", \"85-7=\", \"36+19=\", \"57+4=\", \"42-29=\"),\n    @(\"60-44=\", \"46+47=\", \"53-16=\", \"42-33=\", \"41-39=\"),\n    @(\"16+6=\", \"93-25=\", \"43+29=\", \"52-49=\", \"35-18=\"),\n    @(\"49+33=\", \"87+4=\", \"32-5=\", \"21-3=\", \"62-25=\"),\n    @(\"87-59=\", \"16+78=\", \"13+8=\", \"9+19=\", \"37+14=\"),\n    @(\"62-53=\", \"56+35=\", \"40-21=\", \"86-69=\", \"84-77=\"),\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rowCount = [Math]::Min($tbl.Rows.Count, $newValues.Count)\nfor ($r = 0; $r -lt $rowCount; $r++) {\n    $rowValues = $newValues[$r]\n    $colCount = [Math]::Min($tbl.Rows.Item($r + 1).Cells.Count, $rowValues.Count)\n    for ($c = 0; $c -lt $colCount; $c++) {\n        $cell = $tbl.Cell($r + 1, $c + 1)\n        $cell.Range.Text = $rowValues[$c]\n    }\n}\n"}
